$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Rows 2,3,5,6,7,8 get cycled: each row's record data is replaced by the
# record that used to sit in a different row. Only the columns that actually
# differ between the source/target rows are touched.
#   new row 2 <- old row 8
#   new row 8 <- old row 2
#   new row 3 <- old row 7
#   new row 7 <- old row 6
#   new row 6 <- old row 5
#   new row 5 <- old row 3
# ---------------------------------------------------------------------------

# new row 2 <- old row 8
$ws.Cells.Item(2,1).Value2  = 111749006
$ws.Cells.Item(2,2).Value2  = 8377
$ws.Cells.Item(2,4).Value2  = "LC"
$ws.Cells.Item(2,5).Value2  = 106545
$ws.Cells.Item(2,6).Value2  = "Mindre märgborre"
$ws.Cells.Item(2,7).Value2  = "Tomicus minor"
$ws.Cells.Item(2,8).Value2  = "(Hartig, 1834)"
$ws.Cells.Item(2,17).Value2 = 575512.2089522779
$ws.Cells.Item(2,18).Value2 = 6509825.662577543

# new row 8 <- old row 2
$ws.Cells.Item(8,1).Value2  = 111747186
$ws.Cells.Item(8,2).Value2  = 78107
$ws.Cells.Item(8,4).Value2  = "NT"
$ws.Cells.Item(8,5).Value2  = 6453
$ws.Cells.Item(8,6).Value2  = "Vedskivlav"
$ws.Cells.Item(8,7).Value2  = "Hertelidea botryosa"
$ws.Cells.Item(8,8).Value2  = "(Fr.) Printzen & Kantvilas"
$ws.Cells.Item(8,17).Value2 = 575435.6246570286
$ws.Cells.Item(8,18).Value2 = 6509856.898648335

# new row 3 <- old row 7 (only Id / Ost / Nord differ)
$ws.Cells.Item(3,1).Value2  = 111749883
$ws.Cells.Item(3,17).Value2 = 575336.5075504743
$ws.Cells.Item(3,18).Value2 = 6509789.003789719

# new row 7 <- old row 6
$ws.Cells.Item(7,1).Value2  = 111749097
$ws.Cells.Item(7,2).Value2  = 93388
$ws.Cells.Item(7,4).Value2  = "LC"
$ws.Cells.Item(7,5).Value2  = 2180
$ws.Cells.Item(7,6).Value2  = "Blåmossa"
$ws.Cells.Item(7,7).Value2  = "Leucobryum glaucum"
$ws.Cells.Item(7,8).Value2  = "(Hedw.) Ångstr."
$ws.Cells.Item(7,16).Value2 = "Lilla gruvan (Lilla gruvan), Ög"
$ws.Cells.Item(7,17).Value2 = 575501.7342092508
$ws.Cells.Item(7,18).Value2 = 6509775.591426332
$ws.Cells.Item(7,19).Value2 = 3

# new row 6 <- old row 5
$ws.Cells.Item(6,1).Value2  = 111749860
$ws.Cells.Item(6,2).Value2  = 78107
$ws.Cells.Item(6,4).Value2  = "NT"
$ws.Cells.Item(6,5).Value2  = 6453
$ws.Cells.Item(6,6).Value2  = "Vedskivlav"
$ws.Cells.Item(6,7).Value2  = "Hertelidea botryosa"
$ws.Cells.Item(6,8).Value2  = "(Fr.) Printzen & Kantvilas"
$ws.Cells.Item(6,16).Value2 = "Älgsjöhåll (Älgsjöhåll), Ög"
$ws.Cells.Item(6,17).Value2 = 575356.6078101217
$ws.Cells.Item(6,18).Value2 = 6509772.251964441
$ws.Cells.Item(6,19).Value2 = 1

# new row 5 <- old row 3 (only Id / Ost / Nord differ)
$ws.Cells.Item(5,1).Value2  = 111749897
$ws.Cells.Item(5,17).Value2 = 575336.6687912485
$ws.Cells.Item(5,18).Value2 = 6509780.695668718

# ---------------------------------------------------------------------------
# Rows 10 and 12 are fully swapped (whole record exchanged between the two
# rows), including the presence/absence of a few sparsely-populated columns
# (I, J, M, AF).
# ---------------------------------------------------------------------------

# new row 10 <- old row 12
# I10 must end up as a *text* "25" (not a number), and AF10 must exist as an
# empty placeholder cell - copy those two straight from row 12's current
# (pre-swap) contents so the stored cell type matches exactly. Likewise stash
# row 10's own (empty-placeholder) I-cell so it can be copied into row 12
# afterwards instead of just clearing it away.
$ws.Cells.Item(10,9).Copy($ws.Cells.Item(99,9))
$ws.Cells.Item(12,9).Copy($ws.Cells.Item(10,9))
$ws.Cells.Item(12,32).Copy($ws.Cells.Item(10,32))

$ws.Cells.Item(10,1).Value2  = 111964550
$ws.Cells.Item(10,2).Value2  = 103288
$ws.Cells.Item(10,4).Value2  = "LC"
$ws.Cells.Item(10,5).Value2  = 221144
$ws.Cells.Item(10,6).Value2  = "Grönpyrola"
$ws.Cells.Item(10,7).Value2  = "Pyrola chlorantha"
$ws.Cells.Item(10,8).Value2  = "Sw."
$ws.Cells.Item(10,10).Value2 = "plantor/tuvor"
$ws.Cells.Item(10,11).Value2 = "överblommad"
$ws.Cells.Item(10,13).ClearContents()
$ws.Cells.Item(10,26).Value2 = "00:00"
$ws.Cells.Item(10,28).Value2 = "00:00"

# new row 12 <- old row 10
$ws.Cells.Item(12,1).Value2  = 111964494
$ws.Cells.Item(12,2).Value2  = 56414
$ws.Cells.Item(12,4).Value2  = "NT"
$ws.Cells.Item(12,5).Value2  = 100049
$ws.Cells.Item(12,6).Value2  = "Spillkråka"
$ws.Cells.Item(12,7).Value2  = "Dryocopus martius"
$ws.Cells.Item(12,8).Value2  = "(Linnaeus, 1758)"
$ws.Cells.Item(99,9).Copy($ws.Cells.Item(12,9))
$ws.Cells.Item(99,9).ClearContents()
$ws.Cells.Item(12,10).ClearContents()
$ws.Cells.Item(12,11).Value2 = "adult"
$ws.Cells.Item(12,13).Value2 = "förbiflygande"
$ws.Cells.Item(12,26).Value2 = "10:30"
$ws.Cells.Item(12,28).Value2 = "10:30"
$ws.Cells.Item(12,32).ClearContents()
